# Auto-generated edit script applying numeric value changes per the commit diff
# to the 'Ultros_Profits' sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2676.875
$ws.Range("I2").Value = 226.4
$ws.Range("J2").Value = 6761
$ws.Range("K2").Value = 226.4
$ws.Range("L2").Value = 6761
$ws.Range("M2").Value = -113.4
$ws.Range("N2").Value = -6987
$ws.Range("H15").Value = 1326
$ws.Range("I15").Value = 1326
$ws.Range("K15").Value = 3978
$ws.Range("M15").Value = -3809
$ws.Range("H17").Value = 7571.8
$ws.Range("J17").Value = 5699
$ws.Range("L17").Value = 17097
$ws.Range("N17").Value = -17433
$ws.Range("H138").Value = 2908.8167
$ws.Range("J138").Value = 4818.276
$ws.Range("L138").Value = 14454.828
$ws.Range("N138").Value = -24734.828

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 16182.477
$ws.Range("I2").Value = 23011.857
$ws.Range("J2").Value = 2523.7144
$ws.Range("K2").Value = 23011.857
$ws.Range("L2").Value = 2523.7144
$ws.Range("M2").Value = -22898.857
$ws.Range("N2").Value = -2749.7144
$ws.Range("H32").Value = 2193.66
$ws.Range("I32").Value = 2104.9678
$ws.Range("K32").Value = 2104.9678
$ws.Range("M32").Value = -1817.9678
$ws.Range("H74").Value = 2209.6667
$ws.Range("I74").Value = 1515.5
$ws.Range("K74").Value = 1515.5
$ws.Range("M74").Value = -641.5
$ws.Range("H77").Value = 2209.6667
$ws.Range("I77").Value = 1515.5
$ws.Range("K77").Value = 7577.5
$ws.Range("M77").Value = -3209.5
$ws.Range("H116").Value = 16182.477
$ws.Range("I116").Value = 23011.857
$ws.Range("J116").Value = 2523.7144
$ws.Range("K116").Value = 23011.857
$ws.Range("L116").Value = 2523.7144
$ws.Range("M116").Value = -20717.857
$ws.Range("N116").Value = -7111.7144
$ws.Range("H122").Value = 4950.6665
$ws.Range("I122").Value = 3587.3333
$ws.Range("K122").Value = 10761.9999
$ws.Range("M122").Value = -8311.999899999999
$ws.Range("H132").Value = 1587
$ws.Range("I132").Value = 1587
$ws.Range("K132").Value = 4761
$ws.Range("M132").Value = -2231

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 16182.477
$ws.Range("I3").Value = 23011.857
$ws.Range("J3").Value = 2523.7144
$ws.Range("K3").Value = 23011.857
$ws.Range("L3").Value = 2523.7144
$ws.Range("M3").Value = -22897.857
$ws.Range("N3").Value = -2751.7144
$ws.Range("H22").Value = 236.76923
$ws.Range("I22").Value = 250.5
$ws.Range("J22").Value = 72
$ws.Range("K22").Value = 250.5
$ws.Range("L22").Value = 72
$ws.Range("M22").Value = -77.5
$ws.Range("N22").Value = -418
$ws.Range("H105").Value = 4805.2354
$ws.Range("I105").Value = 3744.8333
$ws.Range("K105").Value = 3744.8333
$ws.Range("M105").Value = -1997.8333

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6171.25
$ws.Range("I16").Value = 6774.4
$ws.Range("K16").Value = 6774.4
$ws.Range("M16").Value = -6487.4
$ws.Range("H113").Value = 6171.25
$ws.Range("I113").Value = 6774.4
$ws.Range("K113").Value = 6774.4
$ws.Range("M113").Value = -4604.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 2200
$ws.Range("I45").Value = 2200
$ws.Range("K45").Value = 6600
$ws.Range("M45").Value = -6068
$ws.Range("H92").Value = 677.7778
$ws.Range("I92").Value = 400
$ws.Range("J92").Value = 1025
$ws.Range("K92").Value = 1200
$ws.Range("L92").Value = 3075
$ws.Range("M92").Value = 48
$ws.Range("N92").Value = -5571
$ws.Range("H97").Value = 475
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 475
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 1425
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -2417
$ws.Range("H98").Value = 719.2
$ws.Range("I98").Value = 774
$ws.Range("J98").Value = 500
$ws.Range("K98").Value = 2322
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = -824
$ws.Range("N98").Value = -4496
$ws.Range("H109").Value = 833.8
$ws.Range("I109").Value = 750.5
$ws.Range("K109").Value = 2251.5
$ws.Range("M109").Value = -1211.5
$ws.Range("H131").Value = 835804.6
$ws.Range("I131").Value = 1112350.8
$ws.Range("J131").Value = 6166.3335
$ws.Range("K131").Value = 3337052.4
$ws.Range("L131").Value = 18499.0005
$ws.Range("M131").Value = -3332012.4
$ws.Range("N131").Value = -28579.0005
$ws.Range("H134").Value = 3650.1538
$ws.Range("I134").Value = 993.1429000000001
$ws.Range("J134").Value = 6750
$ws.Range("K134").Value = 2979.4287
$ws.Range("L134").Value = 20250
$ws.Range("M134").Value = 2090.5713
$ws.Range("N134").Value = -30390
$ws.Range("H139").Value = 250006750

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H80").Value = 20905410
$ws.Range("I80").Value = 187518.33
$ws.Range("J80").Value = 33336144
$ws.Range("K80").Value = 187518.33
$ws.Range("L80").Value = 33336144
$ws.Range("M80").Value = -186520.33
$ws.Range("N80").Value = -33338140
$ws.Range("H83").Value = 20905410
$ws.Range("I83").Value = 187518.33
$ws.Range("J83").Value = 33336144
$ws.Range("K83").Value = 937591.6499999999
$ws.Range("L83").Value = 166680720
$ws.Range("M83").Value = -932599.6499999999
$ws.Range("N83").Value = -166690704
$ws.Range("H92").Value = 15090.2
$ws.Range("J92").Value = 15090.2
$ws.Range("L92").Value = 15090.2
$ws.Range("N92").Value = -18834.2
$ws.Range("H132").Value = 2956.6365
$ws.Range("J132").Value = 2940
$ws.Range("L132").Value = 8820
$ws.Range("N132").Value = -13880

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4878.657
$ws.Range("I132").Value = 3005.2273
$ws.Range("K132").Value = 9015.6819
$ws.Range("M132").Value = -6485.6819
$ws.Range("H136").Value = 3595.8293
$ws.Range("I136").Value = 3261.9167
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 9785.750100000001
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -7235.750100000001
$ws.Range("N136").Value = -23100

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 20000
$ws.Range("I15").Value = 20000
$ws.Range("K15").Value = 20000
$ws.Range("M15").Value = -19712
$ws.Range("H16").Value = 56995
$ws.Range("J16").Value = 56995
$ws.Range("L16").Value = 56995
$ws.Range("N16").Value = -57579
$ws.Range("H46").Value = 64000
$ws.Range("J46").Value = 64000
$ws.Range("L46").Value = 64000
$ws.Range("N46").Value = -64462
$ws.Range("H107").Value = 550.5769
$ws.Range("I107").Value = 604.4761999999999
$ws.Range("K107").Value = 1813.4286
$ws.Range("M107").Value = 106.5714000000003
$ws.Range("H113").Value = 878.86957
$ws.Range("I113").Value = 518.58826
$ws.Range("K113").Value = 1555.76478
$ws.Range("M113").Value = 614.23522
$ws.Range("H122").Value = 2680.853
$ws.Range("I122").Value = 2715.4583
$ws.Range("J122").Value = 2597.8
$ws.Range("K122").Value = 8146.374899999999
$ws.Range("L122").Value = 7793.400000000001
$ws.Range("M122").Value = -5696.374899999999
$ws.Range("N122").Value = -12693.4
$ws.Range("H132").Value = 2387.9707
$ws.Range("I132").Value = 1627.8889
$ws.Range("J132").Value = 5319.7144
$ws.Range("K132").Value = 4883.6667
$ws.Range("L132").Value = 15959.1432
$ws.Range("M132").Value = -2353.6667
$ws.Range("N132").Value = -21019.1432
$ws.Range("H134").Value = 64000
$ws.Range("J134").Value = 64000
$ws.Range("L134").Value = 192000
$ws.Range("M134").Value = -197070

